$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,("adana", "ceyhan", "emlakjet", 10, 5679900)
    ,("adana", "kozan", "emlakjet", 10, 3211000)
    ,("adiyaman", "kahta", "emlakjet", 10, 3640500)
    ,("adana", "aladağ", "emlakjet", 10, 3028000)
    ,("adana", "feke", "emlakjet", 10, 2832500)
    ,("afyonkarahisar", "ihsaniye", "emlakjet", 10, 3965000)
    ,("adana", "feke", "emlakjet", 3, 2050000)
    ,("adana", "aladağ", "emlakjet", 1, 5150000)
    ,("afyonkarahisar", "iscehisar", "emlakjet", 1, 1150000)
    ,("şırnak", "uludere", "emlakjet", 1, 600000)
    ,("adana", "ceyhan", "hepsiemlak", 10, 5335000)
    ,("istanbul", "kadıköy", "hepsiemlak", 10, 20910000)
    ,("adana", "ceyhan", "hepsiemlak", 10, 5335000)
    ,("istanbul", "üsküdar", "hepsiemlak", 10, 13387500)
    ,("adana", "ceyhan", "hepsiemlak", 10, 4214000)
    ,("istanbul", "kartal", "hepsiemlak", 10, 5980000)
    ,("adana", "kozan", "hepsiemlak", 10, 3051500)
    ,("adana", "ceyhan", "hepsiemlak", 10, 5670000)
    ,("istanbul", "kadıköy", "hepsiemlak", 10, 20590000)
    ,("adana", "ceyhan", "hepsiemlak", 10, 5615000)
    ,("adana", "ceyhan", "hepsiemlak", 10, 5555000)
    ,("adana", "ceyhan", "hepsiemlak", 10, 4734000)
    ,("istanbul", "ümraniye", "hepsiemlak", 10, 14121900)
    ,("adana", "aladağ", "hepsiemlak", 10, 7036500)
    ,("adana", "feke", "hepsiemlak", 1, 5500000)
    ,("adana", "feke", "hepsiemlak", 1, 5500000)
    ,("adana", "ceyhan", "hepsiemlak", 10, 4650000)
)

$startRow = 41
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
